# Applies the "Updated cryptos list" data refresh described by the commit.
# Rows 2-51 hold one inline-string row each: A=rank(unchanged), B=Coin,
# C=Link, D=Price, E=Volume(1h). This script only rewrites the D/E (and, for
# the Hedera/InjectiveProtocol swap, B/C) cells that actually changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Price cells such as "7.58" or "0.0000118" parse as real numbers, so a
    # plain .Value assignment would silently turn them into floats/sci-notation.
    # Force the cell to Text, assign the literal string, then restore the
    # default ("Normal") style so no stray number format is left behind -
    # matching the plain inlineStr cells (no "s" attribute) in the target file.
    $cell = $ws.Range($cellRef)
    $looksNumeric = $text -match '^-?\d+(\.\d+)?$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue "D2" "68.032.46"
Set-TextValue "E2" "  +0.34%  "
Set-TextValue "D3" "3.242.88"
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "581.36"
Set-TextValue "E5" "  -0.45%  "
Set-TextValue "D6" "184.74"
Set-TextValue "E6" "  +1.17%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "0.601"
Set-TextValue "E8" "  +0.76%  "
Set-TextValue "E9" "  -3.27%  "
Set-TextValue "E10" "  -0.99%  "
Set-TextValue "E11" "  +0.11%  "
Set-TextValue "D12" "3.810.39"
Set-TextValue "E12" "  +0.03%  "
Set-TextValue "E13" "  +0.14%  "
Set-TextValue "D14" "27.81"
Set-TextValue "E14" "  -2.87%  "
Set-TextValue "D15" "68.042.29"
Set-TextValue "E15" "  +0.37%  "
Set-TextValue "E16" "  -0.87%  "
Set-TextValue "D17" "3.242.96"
Set-TextValue "E17" "  -0.28%  "
Set-TextValue "E18" "  -0.52%  "
Set-TextValue "D19" "13.46"
Set-TextValue "E19" "  -0.60%  "
Set-TextValue "D20" "396.61"
Set-TextValue "E20" "  +4.42%  "
Set-TextValue "D21" "7.58"
Set-TextValue "E21" "  -0.65%  "
Set-TextValue "E22" "  +0.10%  "
Set-TextValue "D23" "71.35"
Set-TextValue "E23" "  +0.15%  "
Set-TextValue "D24" "0.515"
Set-TextValue "E24" "  +0.62%  "
Set-TextValue "D25" "0.0000118"
Set-TextValue "E25" "  -0.88%  "
Set-TextValue "D26" "0.187"
Set-TextValue "E26" "  +2.82%  "
Set-TextValue "D27" "9.63"
Set-TextValue "E27" "  -2.93%  "
Set-TextValue "E28" "  -0.07%  "
Set-TextValue "D29" "1.96"
Set-TextValue "E29" "  -0.86%  "
Set-TextValue "E30" "  -1.17%  "
Set-TextValue "D31" "22.78"
Set-TextValue "E31" "  -0.28%  "
Set-TextValue "E32" "  -0.71%  "
Set-TextValue "E33" "  +0.17%  "
Set-TextValue "E34" "  +0.03%  "
Set-TextValue "D35" "161.87"
Set-TextValue "E35" "  -0.15%  "
Set-TextValue "E36" "  -4.01%  "
Set-TextValue "E37" "  +3.37%  "
Set-TextValue "D38" "26.63"
Set-TextValue "E38" "  +0.73%  "
Set-TextValue "D39" "0.811"
Set-TextValue "E39" "  -2.85%  "
Set-TextValue "D40" "4.59"
Set-TextValue "E40" "  +0.36%  "
Set-TextValue "E41" "  -3.19%  "
Set-TextValue "D42" "2.48"
Set-TextValue "E42" "  -3.74%  "
Set-TextValue "D43" "41.21"
Set-TextValue "E43" "  +0.00%  "
Set-TextValue "B44" "Hedera"
Set-TextValue "C44" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D44" "0.0685"
Set-TextValue "E44" "  -0.04%  "
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D45" "25.20"
Set-TextValue "E45" "  -0.88%  "
Set-TextValue "D46" "2.610.49"
Set-TextValue "E46" "  -0.33%  "
Set-TextValue "D47" "336.24"
Set-TextValue "E47" "  -2.75%  "
Set-TextValue "E48" "  -1.42%  "
Set-TextValue "E49" "  +2.06%  "
Set-TextValue "E50" "  -1.10%  "
Set-TextValue "E51" "  +2.54%  "
